$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column Q (27-jun)
$ws.Range("Q1").Value = "27-jun"

# Fill in Q2:Q18 values (row order matches A2:A18)
$values = @(
    0,
    14.815379981990016,
    14.630213722631952,
    16.255628058144431,
    0,
    6.8780659362289978,
    5.7936075396684261,
    14.858127856878696,
    16.740544025206376,
    12.427535297661905,
    0,
    11.696535955512893,
    0,
    0,
    13.452883396260123,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $values[$i]
}

# Update selection to match the final state
$ws.Range("O7").Select()
